$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2026-01-21 18:50:20"

for ($row = 2; $row -le 12; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
